$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the header formatting (bold font + border) from an existing header
# cell onto the three new header cells before setting their text, so the
# new headers (AD1:AF1) match the style of the rest of row 1 (style index 1).
$ws.Range("A1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122) # xlPasteFormats

# New header labels for the team record columns.
$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Every player row (2-40) gets the same team record values.
for ($r = 2; $r -le 40; $r++) {
    $ws.Cells.Item($r, 30).Value = 70   # AD: Wins
    $ws.Cells.Item($r, 31).Value = 92   # AE: Losses
    $ws.Cells.Item($r, 32).Value = 0    # AF: Ties
}
